{"js": "const replacements = [\n  [\"27\u00d743=\", \"55\u00d739=\"],\n  [\"49\u00d757=\", \"31\u00d796=\"],\n  [\"98\u00d772=\", \"52\u00d788=\"],\n  [\"53\u00d745=\", \"85\u00d764=\"],\n  [\"80\u00d755=\", \"39\u00d746=\"],\n  [\"60\u00d784=\", \"66\u00d783=\"],\n  [\"89\u00d748=\", \"70\u00d775=\"],\n  [\"33\u00d799=\", \"40\u00d773=\"],\n  [\"65\u00d727=\", \"91\u00d712=\"],\n  [\"96\u00d721=\", \"89\u00d774=\"],\n  [\"53\u00d777=\", \"24\u00d767=\"],\n  [\"13\u00d772=\", \"58\u00d772=\"],\n  [\"14\u00d775=\", \"52\u00d734=\"],\n  [\"56\u00d756=\", \"79\u00d727=\"],\n  [\"46\u00d726=\", \"18\u00d750=\"],\n  [\"28\u00d738=\", \"92\u00d721=\"],\n  [\"52\u00d737=\", \"69\u00d716=\"],\n  [\"30\u00d779=\", \"92\u00d741=\"],\n  [\"98\u00d796=\", \"40\u00d773=\"],\n  [\"53\u00d739=\", \"46\u00d777=\"],\n  [\"99\u00d787=\", \"14\u00d781=\"],\n  [\"53\u00d772=\", \"20\u00d712=\"],\n  [\"39\u00d788=\", \"74\u00d797=\"],\n  [\"98\u00d739=\", \"99\u00d794=\"],\n  [\"91\u00d787=\", \"40\u00d754=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"27\u00d743=\", \"55\u00d739=\"),\n    @(\"49\u00d757=\", \"31\u00d796=\"),\n    @(\"98\u00d772=\", \"52\u00d788=\"),\n    @(\"53\u00d745=\", \"85\u00d764=\"),\n    @(\"80\u00d755=\", \"39\u00d746=\"),\n    @(\"60\u00d784=\", \"66\u00d783=\"),\n    @(\"89\u00d748=\", \"70\u00d775=\"),\n    @(\"33\u00d799=\", \"40\u00d773=\"),\n    @(\"65\u00d727=\", \"91\u00d712=\"),\n    @(\"96\u00d721=\", \"89\u00d774=\"),\n    @(\"53\u00d777=\", \"24\u00d767=\"),\n    @(\"13\u00d772=\", \"58\u00d772=\"),\n    @(\"14\u00d775=\", \"52\u00d734=\"),\n    @(\"56\u00d756=\", \"79\u00d727=\"),\n    @(\"46\u00d726=\", \"18\u00d750=\"),\n    @(\"28\u00d738=\", \"92\u00d721=\"),\n    @(\"52\u00d737=\", \"69\u00d716=\"),\n    @(\"30\u00d779=\", \"92\u00d741=\"),\n    @(\"98\u00d796=\", \"40\u00d773=\"),\n    @(\"53\u00d739=\", \"46\u00d777=\"),\n    @(\"99\u00d787=\", \"14\u00d781=\"),\n    @(\"53\u00d772=\", \"20\u00d712=\"),\n    @(\"39\u00d788=\", \"74\u00d797=\"),\n    @(\"98\u00d739=\", \"99\u00d794=\"),\n    @(\"91\u00d787=\", \"40\u00d754=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}"}
